$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 29780
$ws.Range("B2").Value = "Isis da Luz"
$ws.Range("C2").Value = "Marketing"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45093
$ws.Range("G2").Value = 2597.22

# Row 3
$ws.Range("A3").Value = 89503
$ws.Range("B3").Value = "Manuela Aragão"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45103
$ws.Range("G3").Value = 10500.19

# Row 4
$ws.Range("A4").Value = 96106
$ws.Range("B4").Value = "Luiz Otávio da Paz"
$ws.Range("C4").Value = "Financeiro"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45090
$ws.Range("G4").Value = 7802.94

# Row 5
$ws.Range("A5").Value = 85469
$ws.Range("B5").Value = "André Ramos"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("F5").Value = 45093
$ws.Range("G5").Value = 4642.11

# Row 6
$ws.Range("A6").Value = 49419
$ws.Range("B6").Value = "Luiz Gustavo Vieira"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45098
$ws.Range("G6").Value = 3666.54

# Row 7
$ws.Range("A7").Value = 36453
$ws.Range("B7").Value = "Rodrigo Barbosa"
$ws.Range("C7").Value = "Jurídico"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45089
$ws.Range("G7").Value = 5995.42

# Row 8
$ws.Range("A8").Value = 32424
$ws.Range("B8").Value = "Clarice Gonçalves"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("F8").Value = 45083
$ws.Range("G8").Value = 5625.15

# Row 9
$ws.Range("A9").Value = 30437
$ws.Range("B9").Value = "Marcela Rocha"
$ws.Range("C9").Value = "Recursos Humanos"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45089
$ws.Range("G9").Value = 4693.69

# Row 10
$ws.Range("A10").Value = 24513
$ws.Range("B10").Value = "Laís Almeida"
$ws.Range("C10").Value = "Vendas"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 8
$ws.Range("G10").Value = 7835.81

# Row 11
$ws.Range("A11").Value = 69410
$ws.Range("B11").Value = "Luiz Felipe Oliveira"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45091
$ws.Range("G11").Value = 10033.98
